$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A2 to be a numeric value instead of text
$ws.Range("A2").Value = 79174445

# Add new row 3 for the redemption record
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "79174445"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = "2025-08-18T08:51:16"
